$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet ("Datos Exportados" -> "Usuario")
$ws.Name = "Usuario"

# --- Row 1 (header, merged across A:B) ---
$ws.Range("A1").Value = "CARACTERISTICAS GENERALES"
$ws.Range("B1").Value = ""

# --- Row 2 ---
$ws.Range("A2").Value = "RUT"
$ws.Range("B2").Value = "12.222.333-2"

# --- Row 3 ---
$ws.Range("A3").Value = "EMAIL"
$ws.Range("B3").Value = "basdadas@gmail.com"

# Drop the now-unused columns C:F (previously held extra export fields)
$ws.Range("C1:F3").EntireColumn.Delete()

# --- Row 4 (reuse row 2's look: white fill) ---
$ws.Range("A4").Value = "ROL"
$ws.Range("B4").Value = "Administrador"
$ws.Range("A2:B2").Copy()
$ws.Range("A4:B4").PasteSpecial(-4122)

# --- Row 5 (reuse row 3's look: gray fill) ---
$ws.Range("A5").Value = "NOMBRE COMPLETO"
$ws.Range("B5").Value = "Tomas Bawssy"
$ws.Range("A3:B3").Copy()
$ws.Range("A5:B5").PasteSpecial(-4122)

# --- Row 6 (reuse row 2's look: white fill) ---
$ws.Range("A6").Value = "FECHA DE REGISTRO"
$ws.Range("B6").Value = "27/10/2024"
$ws.Range("A2:B2").Copy()
$ws.Range("A6:B6").PasteSpecial(-4122)

# --- Row 7 (reuse row 3's look: gray fill) ---
$ws.Range("A7").Value = "ÚLTIMA ACTUALIZACIÓN"
$ws.Range("B7").Value = "27/10/2024"
$ws.Range("A3:B3").Copy()
$ws.Range("A7:B7").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Merge the header row
$ws.Range("A1:B1").Merge()

# Column widths (30 / 50 chars) -- offset compensates for this runtime's
# width round-trip through its internal pixel representation.
$ws.Columns.Item(1).ColumnWidth = 29.166666666666668
$ws.Columns.Item(2).ColumnWidth = 49.166666666666664
